# Word COM-interop script applying the Assignment_1.docx edit:
#   1. "...represented through symbol i or j. " -> "...represented through symbol j. "
#   2. "Complex data type example: X = 10 + 5*i" -> "Complex data type example: X = 10 + 5*j"
#   3. Remove the "age":38, field from the D1 dictionary example text.

$d = $word.ActiveDocument

# 1) Drop the "i or" before "j." in the numeric-types explanation paragraph.
$d.Content.Find.Execute(
    "symbol i or j.", $true, $false, $false, $false, $false,
    $true, 1, $false, "symbol j.", 2) | Out-Null

# 2) Change the complex-number example from 5*i to 5*j.
$d.Content.Find.Execute(
    "X = 10 + 5*i", $true, $false, $false, $false, $false,
    $true, 1, $false, "X = 10 + 5*j", 2) | Out-Null

# 3) Remove the age field from the dictionary literal example.
$d.Content.Find.Execute(
    [char]8221 + ", " + [char]8220 + "age" + [char]8221 + ":38, " + [char]8220,
    $true, $false, $false, $false, $false,
    $true, 1, $false, [char]8221 + ", " + [char]8220, 2) | Out-Null

Write-Output "done"
